$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the column headers: "<name>_old" -> "<name>_FV2210"
#    and "<name>_new" -> "<name>_FV2304" (the "diff" header in K1 is left
#    untouched).
# ---------------------------------------------------------------------------
$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")
$oldCols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$newCols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Range($oldCols[$i] + "1").Value = $baseNames[$i] + "_FV2210"
    $ws.Range($newCols[$i] + "1").Value = $baseNames[$i] + "_FV2304"
}

# ---------------------------------------------------------------------------
# 2) Turn the data range A1:U52 into an actual Excel Table ("Table1") so the
#    header row doubles as a filterable table header. We preserve the header
#    row's pre-existing formatting (bold / shaded / bordered style) rather
#    than letting Excel bake it into a brand-new headerRowDxfId: stash a copy
#    of the header formatting on a scratch row, strip the header's direct
#    formatting before creating the table, create the table, then paste the
#    original formatting back onto the header and discard the scratch copy.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("A100:U100")
$headerRange.Copy($scratch)

$headerRange.ClearFormats()

$dataRange = $ws.Range("A1:U52")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

$scratch.Copy()
$headerRange.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$scratch.Clear()

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split/freeze below row 1).
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
